$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.484.92'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '3.559.61'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D5').Value = "'611.06"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.03%  '
$ws.Range('D6').Value = "'173.07"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('D8').Value = '3.554.64'
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +2.90%  '
$ws.Range('D11').Value = "'7.55"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +12.49%  '
$ws.Range('D12').Value = "'0.586"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').Value = "'46.61"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').Value = '4.143.09'
$ws.Range('D16').Value = "'8.37"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.27%  '
$ws.Range('D17').Value = "'614.59"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.27%  '
$ws.Range('D18').Value = '3.562.33'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '70.608.01'
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('D21').Value = "'17.39"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = "'9.40"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -16.30%  '
$ws.Range('D24').Value = "'16.11"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = "'97.04"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('D26').Value = "'3.83"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('D29').Value = "'33.47"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').Value = "'9.06"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('D31').Value = "'8.51"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('E32').Value = '  -3.44%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').Value = "'1.30"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'6.97"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D35').Value = "'576.57"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -9.13%  '
$ws.Range('D36').Value = "'3.68"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('D40').Value = "'57.39"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('D43').Value = '3.387.82'
$ws.Range('E44').Value = '  -3.12%  '
$ws.Range('D45').Value = "'33.18"
$ws.Range('D45').ClearFormats()
$ws.Range('E46').Value = '  +7.07%  '
$ws.Range('D47').Value = '0.0₃0705'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').Value = "'2.61"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').Value = "'133.86"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.34%  '
